# Commit: "Fruta / hortaliza, semanal"
# A new weekly price-report row for "Ajo" (Chilote variety) is inserted
# into the data table at row 178, pushing every subsequent row down by
# one (old row 178 becomes 179, ..., old row 273 becomes 274).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 178 (shifts rows 178:273 -> 179:274)
$ws.Rows("178:178").Insert()

# Populate the newly inserted row with the new record
$ws.Range("A178").Value2 = 4
$ws.Range("B178").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C178").Value2 = "Los Lagos"
$ws.Range("D178").Value2 = 44719
$ws.Range("E178").Value2 = 10
$ws.Range("F178").Value2 = 100112003
$ws.Range("G178").Value2 = "Ajo"
$ws.Range("H178").Value2 = "Chilote"
$ws.Range("I178").Value2 = "Primera"
$ws.Range("J178").Value2 = 210
$ws.Range("K178").Value2 = 22000
$ws.Range("L178").Value2 = 22000
$ws.Range("M178").Value2 = 22000
$ws.Range("N178").Value2 = "$/caja 10 kilos"
$ws.Range("O178").Value2 = "China"
$ws.Range("P178").Value2 = 2200
$ws.Range("Q178").Value2 = 10
$ws.Range("R178").Value2 = "Hortaliza"
